$wb = $excel.ActiveWorkbook

# --- Update the Date property on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-02-05T08:09:31+00:00"

# --- Remove the "valeur" (Value of the observation) row from the Elements sheet ---
# It currently lives in row 7; deleting it shifts the following "evaluation" row
# (previously row 8) up to become the new row 7, and drops the now-unused
# shared strings ("fr-lm-group-de-questionnaires-devaluation.valeur" and
# "Valeur de l'observation") automatically.
$elements = $wb.Worksheets.Item("Elements")
$elements.Rows("7").Delete()
